$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.447159934046681
$ws.Range("C2").Value = 0.1183654977525634
$ws.Range("D2").Value = 0.05647233328883239
$ws.Range("E2").Value = 0.1211252422109226
$ws.Range("F2").Value = 1.450653168606735
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 1.012805003001851
$ws.Range("K2").Value = 0.5539443074843007
$ws.Range("M2").Value = 0.2907387100251739
$ws.Range("N2").Value = 2.175246244320618

$ws.Range("B3").Value = 0.4095024174869195
$ws.Range("C3").Value = 0.107538064383391
$ws.Range("D3").Value = 0.05652549211416336
$ws.Range("E3").Value = 0.1114079124986063
$ws.Range("F3").Value = 1.43216799730682
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 1.008086454538493
$ws.Range("K3").Value = 0.5062130379697578
$ws.Range("M3").Value = 0.2664624142609497
$ws.Range("N3").Value = 2.185885919413813

$ws.Range("B4").Value = 0.3865995802259192
$ws.Range("C4").Value = 0.1009493183574079
$ws.Range("D4").Value = 0.0565567656189927
$ws.Range("E4").Value = 0.1055151441356088
$ws.Range("F4").Value = 1.421597236561709
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 1.005664435391232
$ws.Range("K4").Value = 0.4771790985443545
$ws.Range("M4").Value = 0.2517154666154724
$ws.Range("N4").Value = 2.193064949504105

$ws.Range("B5").Value = 0.3773215009837259
$ws.Range("C5").Value = 0.09827912791156734
$ws.Range("D5").Value = 0.05656917047867616
$ws.Range("E5").Value = 0.1031320722319222
$ws.Range("F5").Value = 1.417485132810938
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 1.004796747489685
$ws.Range("K5").Value = 0.4654160143265642
$ws.Range("M5").Value = 0.2457455522015977
$ws.Range("N5").Value = 2.196152724653402

$ws.Range("B6").Value = 0.3757842039003663
$ws.Range("C6").Value = 0.09783663311733903
$ws.Range("D6").Value = 0.05657120992722753
$ws.Range("E6").Value = 0.1027374615082124
$ws.Range("F6").Value = 1.416814122139485
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 1.004659868909251
$ws.Range("K6").Value = 0.4634668927290022
$ws.Range("M6").Value = 0.244756633405359
$ws.Range("N6").Value = 2.196675240725384

$ws.Range("B7").Value = 0.3864742301406352
$ws.Range("C7").Value = 0.1009132475935104
$ws.Range("D7").Value = 0.05655693428369268
$ws.Range("E7").Value = 0.1054829315198518
$ws.Range("F7").Value = 1.421540987946059
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 1.005652250639777
$ws.Range("K7").Value = 0.4770201804957424
$ws.Range("M7").Value = 0.2516347944289805
$ws.Range("N7").Value = 2.193105935653165

$ws.Range("B8").Value = 0.4341301106053095
$ws.Range("C8").Value = 0.1146197913984395
$ws.Range("D8").Value = 0.05649094849349012
$ws.Range("E8").Value = 0.11775925052509
$ws.Range("F8").Value = 1.444117522820065
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 1.011079306099084
$ws.Range("K8").Value = 0.5374296969102943
$ws.Range("M8").Value = 0.2823350696118325
$ws.Range("N8").Value = 2.178780536495012

$ws.Range("B9").Value = 0.5293293384289086
$ws.Range("C9").Value = 0.1419771266253065
$ws.Range("D9").Value = 0.05635052502541704
$ws.Range("E9").Value = 0.1424307895993522
$ws.Range("F9").Value = 1.494593657150844
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 1.025502279070814
$ws.Range("K9").Value = 0.6580778015864155
$ws.Range("M9").Value = 0.3438170289361224
$ws.Range("N9").Value = 2.155826967318603

$ws.Range("B10").Value = 0.600355990660006
$ws.Range("C10").Value = 0.1623816424973938
$ws.Range("D10").Value = 0.05624038460927494
$ws.Range("E10").Value = 0.1609409972825446
$ws.Range("F10").Value = 1.535495407624083
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 1.038420108519304
$ws.Range("K10").Value = 0.7480834786603339
$ws.Range("M10").Value = 0.3897982833435023
$ws.Range("N10").Value = 2.14211117326974

$ws.Range("B11").Value = 0.6329078513584534
$ws.Range("C11").Value = 0.1717334078409465
$ws.Range("D11").Value = 0.05618871546012016
$ws.Range("E11").Value = 0.1694495231298987
$ws.Range("F11").Value = 1.554939478886979
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 1.04480457637284
$ws.Range("K11").Value = 0.7893337248194996
$ws.Range("M11").Value = 0.410899300822031
$ws.Range("N11").Value = 2.136558583585952

$ws.Range("B12").Value = 0.6452693322060838
$ws.Range("C12").Value = 0.1752848905518647
$ws.Range("D12").Value = 0.05616892107542348
$ws.Range("E12").Value = 0.1726844711147066
$ws.Range("F12").Value = 1.562423431925851
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 1.047295532701725
$ws.Range("K12").Value = 0.8049985816891478
$ws.Range("M12").Value = 0.4189166185275042
$ws.Range("N12").Value = 2.134555017945445

$ws.Range("B13").Value = 0.6426055172898373
$ws.Range("C13").Value = 0.1745195602854608
$ws.Range("D13").Value = 0.05617319435739354
$ws.Range("E13").Value = 0.1719871865815605
$ws.Range("F13").Value = 1.560806244203647
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 1.046755797007052
$ws.Range("K13").Value = 0.8016229008923972
$ws.Range("M13").Value = 0.4171887468803988
$ws.Range("N13").Value = 2.134982109924181

$ws.Range("B14").Value = 0.6339241399979869
$ws.Range("C14").Value = 0.1720253856929332
$ws.Range("D14").Value = 0.05618709155978507
$ws.Range("E14").Value = 0.1697154030222379
$ws.Range("F14").Value = 1.555552762485988
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 1.045008038511874
$ws.Range("K14").Value = 0.7906215930508722
$ws.Range("M14").Value = 0.4115583507014975
$ws.Range("N14").Value = 2.136391761487516

$ws.Range("B15").Value = 0.6286110796554567
$ws.Range("C15").Value = 0.1704989611269809
$ws.Range("D15").Value = 0.05619557416068943
$ws.Range("E15").Value = 0.1683255653134239
$ws.Range("F15").Value = 1.552350612142249
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 1.043947038422985
$ws.Range("K15").Value = 0.7838887489481579
$ws.Range("M15").Value = 0.408113074421486
$ws.Range("N15").Value = 2.13726812641454

$ws.Range("B16").Value = 0.5982335153723
$ws.Range("C16").Value = 0.1617718975122955
$ws.Range("D16").Value = 0.05624372952049583
$ws.Range("E16").Value = 0.1603867425453345
$ws.Range("F16").Value = 1.534241586109545
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 1.038013112491356
$ws.Range("K16").Value = 0.7453938611702426
$ws.Range("M16").Value = 0.3884230149123908
$ws.Range("N16").Value = 2.142487895312996

$ws.Range("B17").Value = 0.5796597633301701
$ws.Range("C17").Value = 0.1564360763944421
$ws.Range("D17").Value = 0.05627286787401253
$ws.Range("E17").Value = 0.1555393047417226
$ws.Range("F17").Value = 1.523347149060911
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 1.034503139705471
$ws.Range("K17").Value = 0.7218570996028859
$ws.Range("M17").Value = 0.3763911707499261
$ws.Range("N17").Value = 2.145866184696743

$ws.Range("B18").Value = 0.5689993570857439
$ws.Range("C18").Value = 0.1533735994216272
$ws.Range("D18").Value = 0.05628948028958369
$ws.Range("E18").Value = 0.1527594659934124
$ws.Range("F18").Value = 1.517159774224609
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 1.032532105045881
$ws.Range("K18").Value = 0.7083481837057946
$ws.Range("M18").Value = 0.36948806054356
$ws.Range("N18").Value = 2.147873912245331

$ws.Range("B19").Value = 0.565393827520694
$ws.Range("C19").Value = 0.1523378169049749
$ws.Range("D19").Value = 0.0562950797905355
$ws.Range("E19").Value = 0.1518196738463971
$ws.Range("F19").Value = 1.515078358710284
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 1.031872950681176
$ws.Range("K19").Value = 0.7037792369341389
$ws.Range("M19").Value = 0.367153743995452
$ws.Range("N19").Value = 2.14856478403442

$ws.Range("B20").Value = 0.5816346185089571
$ws.Range("C20").Value = 0.1570034050241702
$ws.Range("D20").Value = 0.05626978130382732
$ws.Range("E20").Value = 0.1560544638830805
$ws.Range("F20").Value = 1.524498719941747
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 1.0348718325627
$ws.Range("K20").Value = 0.7243596448547294
$ws.Range("M20").Value = 0.3776701891324947
$ws.Range("N20").Value = 2.145499869281352

$ws.Range("B21").Value = 0.6364731259485836
$ws.Range("C21").Value = 0.1727577079043385
$ws.Range("D21").Value = 0.05618301583756846
$ws.Range("E21").Value = 0.170382327117963
$ws.Range("F21").Value = 1.557092552459963
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 1.045519406400693
$ws.Range("K21").Value = 0.793851739961724
$ws.Range("M21").Value = 0.4132114036494698
$ws.Range("N21").Value = 2.135975021092221

$ws.Range("B22").Value = 0.6725161247630069
$ws.Range("C22").Value = 0.1831134204102511
$ws.Range("D22").Value = 0.0561249769485066
$ws.Range("E22").Value = 0.1798220506826524
$ws.Range("F22").Value = 1.579099447477731
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 1.052905553272169
$ws.Range("K22").Value = 0.8395270819162022
$ws.Range("M22").Value = 0.4365961474292988
$ws.Range("N22").Value = 2.130327645275131

$ws.Range("B23").Value = 0.6532607107408239
$ws.Range("C23").Value = 0.1775808976674966
$ws.Range("D23").Value = 0.05615607635113307
$ws.Range("E23").Value = 0.1747768788224633
$ws.Range("F23").Value = 1.567289300948886
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 1.048924247357576
$ws.Range("K23").Value = 0.8151255790032508
$ws.Range("M23").Value = 0.4241008248648228
$ws.Range("N23").Value = 2.133288793090884

$ws.Range("B24").Value = 0.5807417307620142
$ws.Range("C24").Value = 0.1567468996992147
$ws.Range("D24").Value = 0.05627117717588526
$ws.Range("E24").Value = 0.1558215386015362
$ws.Range("F24").Value = 1.523977858088642
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 1.03470500045249
$ws.Range("K24").Value = 0.7232281735194874
$ws.Range("M24").Value = 0.3770919008587015
$ws.Range("N24").Value = 2.145665276436688

$ws.Range("B25").Value = 0.5033865067821637
$ws.Range("C25").Value = 0.1345236235841014
$ws.Range("D25").Value = 0.05638972030302547
$ws.Range("E25").Value = 0.1356904722248871
$ws.Range("F25").Value = 1.48027078216721
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 1.021194050073689
$ws.Range("K25").Value = 0.6252019211410413
$ws.Range("M25").Value = 0.3270444468715255
$ws.Range("N25").Value = 2.16148480625148
